# Update countries & provincias Spain
# - Update the "Datos actualizados" timestamp in A1
# - Update several country rows whose totals changed, which in turn
#   changes the relative ranking (and therefore country name) for a
#   handful of rows that are tied/adjacent in the ranking.
# - Update many numeric Casos/Recuperados/Muertes values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update timestamp text in A1 ----
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 17:04"

# ---- Country name (column A) swaps caused by re-ranking ----
$ws.Range("A87").Value  = "Zambia"
$ws.Range("A88").Value  = "Libano"

$ws.Range("A97").Value  = "Albania"
$ws.Range("A98").Value  = "Haiti"
$ws.Range("A99").Value  = "Finlandia"

$ws.Range("A145").Value = "Jordania"
$ws.Range("A146").Value = "Uruguay"

$ws.Range("A164").Value = "Trinidad yTobago"
$ws.Range("A165").Value = "Crucero"
$ws.Range("A166").Value = "San Marino"

$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# ---- Numeric data updates (B:H = Casos totales, Nuevos casos, Casos
#      activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ----

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
    $ws.Range("H$row").Value = $h
}

Set-Row   4 5709066 8135 3063410 2469111 0 208 176545
Set-Row  16  312659    0  233651   72602 0  76   6406
Set-Row  35   89010  883   57734   29771 0   4   1505
Set-Row  49   56099   68   53119    2953 0   0     27
Set-Row  65   31937  522   21885    9131 0   7    921
Set-Row  87   10372  154    9126     972 0   5    274
Set-Row  88   10347    0    2928    7310 0   0    109
Set-Row  89   10190   28    8857    1069 0   2    264
Set-Row  96    8203   37    7006    1131 0   1     66
Set-Row  97    7967  155    3986    3743 0   4    238
Set-Row  98    7949    0    5337    2416 0   0    196
Set-Row  99    7842   37    7100     408 0   0    334
Set-Row 125    2902    0    2765     126 0   0     11
Set-Row 141    1750   94    1194     537 0   3     19
Set-Row 145    1498   16    1261     226 0   0     11
Set-Row 146    1493    0    1228     225 0   0     40
Set-Row 164     730   44     140     578 0   0     12
Set-Row 165     712    0     651      48 0   0     13
Set-Row 166     704    0     657       5 0   0     42
Set-Row 174     383    1     274     109 0   0      0
Set-Row 213      13    0      13       0 0   0      0
Set-Row 214      13    0      12       0 0   0      1
